$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("profiles")

$ws.Range("A4").Value = "sia-estacio"
$ws.Range("B4").Value = "aluno"
$ws.Range("C4").Value = "Aluno teste"
